$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "M_PL" group header (R1:Y1, merged), mirroring M_%cit (B1:I1) and M_ETR (J1:Q1) ---
$ws.Range("R1:Y1").Merge()
$ws.Range("R1").Value = "M_PL"
$ws.Range("J1").Copy()
$ws.Range("R1:Y1").PasteSpecial(-4122)  # xlPasteFormats

# --- Sub-headers for the new group (same 8 labels as the other two groups) ---
$subHeaders = @("GFA - Sales", "GFA - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", "OECD (20%) - Sales", "OECD (20%) - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")
$cols = @("R", "S", "T", "U", "V", "W", "X", "Y")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value = $subHeaders[$i]
}
$ws.Range("J2:Q2").Copy()
$ws.Range("R2:Y2").PasteSpecial(-4122)  # xlPasteFormats

# --- New "M_PL" data values for rows 4-13 ---
$data = @{
    4  = @(64537961721, 64711639466, 62346466255, 62520144000, 70802529578, 70802529578, 70802529578, 70802529578)
    5  = @(976531986457, 988562844368, 976531986457, 988562844368, 1017055371530, 1017055371530, 1017055371530, 1017055371530)
    6  = @(25762595315, 37793453226, 25762595315, 37793453226, 40333624448, 40333624448, 40333624448, 40333624448)
    7  = @(35956611724, 49146191814, 35956611724, 49982592968, 51692957042, 51692957042, 51692957042, 51692957042)
    8  = @(933471841988, 933471841988, 933471841988, 933471841988, 959424197928, 959424197928, 959424197928, 959424197928)
    9  = @(12457548014, 30286449365, 11068720584, 31122850519, 41023269259, 41023269259, 41023269259, 41023269259)
    10 = @(65619795685, 65619795685, 65619795685, 65619795685, 91572151625, 91572151625, 91572151625, 91572151625)
    11 = @(968765904574, 968939582319, 966574409108, 966748086853, 1007717406299, 1007717406299, 1007717406299, 1007717406299)
    12 = @(100913858271, 101087536016, 98722362805, 98896040550, 139865359996, 139865359996, 139865359996, 139865359996)
    13 = @(867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303, 867852046303)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $values[$i]
    }
}
